# Applies the cell-content updates to Sheet1 as described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "71.306.75"
$ws.Cells.Item(2, 5).Value = "  -0.21%  "
$ws.Cells.Item(3, 4).Value = "3.804.87"
$ws.Cells.Item(3, 5).Value = "  -1.17%  "
$ws.Cells.Item(4, 5).Value = "  +0.03%  "
$cell = $ws.Cells.Item(5, 4)
$cell.NumberFormat = "@"
$cell.Value = "704.72"
$cell.ClearFormats()
$ws.Cells.Item(5, 5).Value = "  -1.55%  "
$cell = $ws.Cells.Item(6, 4)
$cell.NumberFormat = "@"
$cell.Value = "171.22"
$cell.ClearFormats()
$ws.Cells.Item(6, 5).Value = "  -1.14%  "
$ws.Cells.Item(7, 4).Value = "3.804.42"
$ws.Cells.Item(7, 5).Value = "  -1.13%  "
$ws.Cells.Item(9, 5).Value = "  -0.34%  "
$ws.Cells.Item(10, 5).Value = "  -2.24%  "
$ws.Cells.Item(11, 5).Value = "  +1.39%  "
$cell = $ws.Cells.Item(12, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.485"
$cell.ClearFormats()
$ws.Cells.Item(12, 5).Value = "  +5.24%  "
$cell = $ws.Cells.Item(13, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.0000250"
$cell.ClearFormats()
$ws.Cells.Item(13, 5).Value = "  -3.03%  "
$ws.Cells.Item(14, 5).Value = "  -1.36%  "
$ws.Cells.Item(15, 4).Value = "4.446.21"
$ws.Cells.Item(15, 5).Value = "  -0.99%  "
$ws.Cells.Item(16, 4).Value = "3.798.60"
$ws.Cells.Item(16, 5).Value = "  -2.73%  "
$ws.Cells.Item(17, 4).Value = "71.479.85"
$ws.Cells.Item(17, 5).Value = "  +0.23%  "
$cell = $ws.Cells.Item(18, 4)
$cell.NumberFormat = "@"
$cell.Value = "7.22"
$cell.ClearFormats()
$ws.Cells.Item(18, 5).Value = "  -0.37%  "
$cell = $ws.Cells.Item(19, 4)
$cell.NumberFormat = "@"
$cell.Value = "17.48"
$cell.ClearFormats()
$ws.Cells.Item(19, 5).Value = "  +0.26%  "
$ws.Cells.Item(20, 5).Value = "  -0.28%  "
$cell = $ws.Cells.Item(21, 4)
$cell.NumberFormat = "@"
$cell.Value = "515.55"
$cell.ClearFormats()
$ws.Cells.Item(21, 5).Value = "  +3.79%  "
$ws.Cells.Item(22, 5).Value = "  -3.01%  "
$cell = $ws.Cells.Item(23, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.715"
$cell.ClearFormats()
$ws.Cells.Item(23, 5).Value = "  -1.70%  "
$cell = $ws.Cells.Item(24, 4)
$cell.NumberFormat = "@"
$cell.Value = "83.95"
$cell.ClearFormats()
$ws.Cells.Item(24, 5).Value = "  -1.71%  "
$ws.Cells.Item(25, 5).Value = "  -4.15%  "
$cell = $ws.Cells.Item(26, 4)
$cell.NumberFormat = "@"
$cell.Value = "12.62"
$cell.ClearFormats()
$ws.Cells.Item(26, 5).Value = "  +3.69%  "
$ws.Cells.Item(27, 4).Value = "3.948.93"
$ws.Cells.Item(27, 5).Value = "  -1.23%  "
$cell = $ws.Cells.Item(28, 4)
$cell.NumberFormat = "@"
$cell.Value = "10.30"
$cell.ClearFormats()
$ws.Cells.Item(28, 5).Value = "  -3.41%  "
$ws.Cells.Item(29, 5).Value = "  -0.09%  "
$ws.Cells.Item(30, 5).Value = "  -4.39%  "
$cell = $ws.Cells.Item(31, 4)
$cell.NumberFormat = "@"
$cell.Value = "3.01"
$cell.ClearFormats()
$ws.Cells.Item(31, 5).Value = "  -6.72%  "
$ws.Cells.Item(32, 2).Value = "NEARProtocol"
$ws.Cells.Item(32, 3).Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$cell = $ws.Cells.Item(32, 4)
$cell.NumberFormat = "@"
$cell.Value = "7.38"
$cell.ClearFormats()
$ws.Cells.Item(32, 5).Value = "  -1.91%  "
$ws.Cells.Item(33, 2).Value = "ImmutableX"
$ws.Cells.Item(33, 3).Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$cell = $ws.Cells.Item(33, 4)
$cell.NumberFormat = "@"
$cell.Value = "2.24"
$cell.ClearFormats()
$ws.Cells.Item(33, 5).Value = "  -0.46%  "
$cell = $ws.Cells.Item(34, 4)
$cell.NumberFormat = "@"
$cell.Value = "29.14"
$cell.ClearFormats()
$ws.Cells.Item(34, 5).Value = "  -1.29%  "
$ws.Cells.Item(35, 5).Value = "  -2.36%  "
$cell = $ws.Cells.Item(36, 4)
$cell.NumberFormat = "@"
$cell.Value = "9.31"
$cell.ClearFormats()
$ws.Cells.Item(36, 5).Value = "  +0.83%  "
$ws.Cells.Item(37, 4).Value = "3.769.58"
$ws.Cells.Item(37, 5).Value = "  -1.12%  "
$ws.Cells.Item(38, 5).Value = "  +0.02%  "
$cell = $ws.Cells.Item(39, 4)
$cell.NumberFormat = "@"
$cell.Value = "6.52"
$cell.ClearFormats()
$ws.Cells.Item(39, 5).Value = "  +7.73%  "
$ws.Cells.Item(40, 5).Value = "  -2.41%  "
$ws.Cells.Item(41, 5).Value = "  +6.89%  "
$ws.Cells.Item(42, 5).Value = "  -2.39%  "
$cell = $ws.Cells.Item(43, 4)
$cell.NumberFormat = "@"
$cell.Value = "3.23"
$cell.ClearFormats()
$ws.Cells.Item(43, 5).Value = "  -4.22%  "
$ws.Cells.Item(45, 5).Value = "  +0.12%  "
$cell = $ws.Cells.Item(46, 4)
$cell.NumberFormat = "@"
$cell.Value = "167.74"
$cell.ClearFormats()
$ws.Cells.Item(46, 5).Value = "  +2.40%  "
$cell = $ws.Cells.Item(47, 4)
$cell.NumberFormat = "@"
$cell.Value = "50.15"
$cell.ClearFormats()
$ws.Cells.Item(47, 5).Value = "  +2.82%  "
$cell = $ws.Cells.Item(48, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.000306"
$cell.ClearFormats()
$ws.Cells.Item(48, 5).Value = "  -4.45%  "
$cell = $ws.Cells.Item(49, 4)
$cell.NumberFormat = "@"
$cell.Value = "427.26"
$cell.ClearFormats()
$ws.Cells.Item(49, 5).Value = "  +1.72%  "
$ws.Cells.Item(50, 5).Value = "  +0.39%  "
$cell = $ws.Cells.Item(51, 4)
$cell.NumberFormat = "@"
$cell.Value = "8.66"
$cell.ClearFormats()
$ws.Cells.Item(51, 5).Value = "  +0.16%  "
